$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-FooterHeaderPicture($story, $newName) {
    if ($story.Exists) {
        $rng = $story.Range
        for ($i = 1; $i -le $rng.InlineShapes.Count; $i++) {
            $ishp = $rng.InlineShapes.Item($i)
            $shp = $ishp.ConvertToShape()
            $shp.Name = $newName
            $shp.ConvertToInlineShape() | Out-Null
        }
    }
}

# Footer 1 (primary footer): PearsonLogo image1.png -> image2.png (id=1)
Rename-FooterHeaderPicture $sec.Footers.Item(1) "image2.png"

# Footer 2 (first-page footer): PearsonLogo image1.png -> image2.png (id=2)
Rename-FooterHeaderPicture $sec.Footers.Item(2) "image2.png"

# Header 2 (first-page header): BTec_Logo-Orange image2.jpg -> image1.jpg (id=3)
Rename-FooterHeaderPicture $sec.Headers.Item(2) "image1.jpg"
